$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Password for anadkarni@essenmed.com changed while re-running the login flow
$ws.Range("B5").Value = "Bronx@1995"

# Leave the cursor where testing last left off
$ws.Range("C7").Select()
